$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
# Row 53
$ws.Range("H53").Value = 12317.2
$ws.Range("I53").Value = 15219
$ws.Range("J53").Value = 10382.667
$ws.Range("K53").Value = 15219
$ws.Range("L53").Value = 10382.667
$ws.Range("M53").Value = -14582
$ws.Range("N53").Value = -11656.667
# Row 61
$ws.Range("H61").Value = 698.3333
$ws.Range("I61").Value = 698.3333
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2094.9999
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -1922.9999
$ws.Range("N61").ClearContents()
# Row 74
$ws.Range("H74").Value = 75007560
$ws.Range("I74").Value = 136366160
$ws.Range("K74").Value = 136366160
$ws.Range("M74").Value = -136365224
# Row 77
$ws.Range("H77").Value = 75007560
$ws.Range("I77").Value = 136366160
$ws.Range("K77").Value = 681830800
$ws.Range("M77").Value = -681826120
# Row 111
$ws.Range("H111").Value = 8930708
$ws.Range("I111").Value = 41668030
$ws.Range("J111").Value = 2347.7273
$ws.Range("K111").Value = 125004090
$ws.Range("L111").Value = 7043.1819
$ws.Range("M111").Value = -125001023
$ws.Range("N111").Value = -13177.1819

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 35715516
$ws.Range("I2").Value = 1059.5
$ws.Range("K2").Value = 1059.5
$ws.Range("M2").Value = -946.5
# Row 32
$ws.Range("H32").Value = 5564759
$ws.Range("I32").Value = 6067622
$ws.Range("K32").Value = 6067622
$ws.Range("M32").Value = -6067335
# Row 37
$ws.Range("H37").Value = 34
$ws.Range("I37").Value = 34
$ws.Range("K37").Value = 34
$ws.Range("M37").Value = 239
# Row 74
$ws.Range("H74").Value = 42888.12
$ws.Range("I74").Value = 60573.176
$ws.Range("J74").Value = 5307.375
$ws.Range("K74").Value = 60573.176
$ws.Range("L74").Value = 5307.375
$ws.Range("M74").Value = -59699.176
$ws.Range("N74").Value = -7055.375
# Row 77
$ws.Range("H77").Value = 42888.12
$ws.Range("I77").Value = 60573.176
$ws.Range("J77").Value = 5307.375
$ws.Range("K77").Value = 302865.88
$ws.Range("L77").Value = 26536.875
$ws.Range("M77").Value = -298497.88
$ws.Range("N77").Value = -35272.875
# Row 102
$ws.Range("H102").Value = 2436.7354
$ws.Range("I102").Value = 1913.0435
$ws.Range("K102").Value = 1913.0435
$ws.Range("M102").Value = -291.0435
# Row 116
$ws.Range("H116").Value = 35715516
$ws.Range("I116").Value = 1059.5
$ws.Range("K116").Value = 1059.5
$ws.Range("M116").Value = 1234.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 35715516
$ws.Range("I3").Value = 1059.5
$ws.Range("K3").Value = 1059.5
$ws.Range("M3").Value = -945.5
# Row 92
$ws.Range("H92").Value = 39999
$ws.Range("J92").Value = 39999
$ws.Range("L92").Value = 39999
$ws.Range("N92").Value = -44991
# Row 105
$ws.Range("H105").Value = 2396.875
$ws.Range("I105").Value = 1674.05
$ws.Range("K105").Value = 1674.05
$ws.Range("M105").Value = 72.95000000000005

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 4879.905
$ws.Range("I62").Value = 4686.125
$ws.Range("J62").Value = 5500
$ws.Range("K62").Value = 4686.125
$ws.Range("L62").Value = 5500
$ws.Range("M62").Value = -4062.125
$ws.Range("N62").Value = -6748
# Row 65
$ws.Range("H65").Value = 4879.905
$ws.Range("I65").Value = 4686.125
$ws.Range("J65").Value = 5500
$ws.Range("K65").Value = 23430.625
$ws.Range("L65").Value = 27500
$ws.Range("M65").Value = -20310.625
$ws.Range("N65").Value = -33740
# Row 86
$ws.Range("H86").Value = 28425568
$ws.Range("J86").Value = 10250
$ws.Range("L86").Value = 10250
$ws.Range("N86").Value = -12496
# Row 89
$ws.Range("H89").Value = 28425568
$ws.Range("J89").Value = 10250
$ws.Range("L89").Value = 51250
$ws.Range("N89").Value = -62482
# Row 132
$ws.Range("H132").Value = 3323.58
$ws.Range("I132").Value = 2484.5278
$ws.Range("K132").Value = 7453.5834
$ws.Range("M132").Value = -4923.5834

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
# Row 12
$ws.Range("H12").Value = 3125507.8
$ws.Range("J12").Value = 4166727
$ws.Range("L12").Value = 12500181
$ws.Range("N12").Value = -12500527
# Row 23
$ws.Range("H23").Value = 541
$ws.Range("J23").Value = 568.625
$ws.Range("L23").Value = 1705.875
$ws.Range("N23").Value = -2175.875
# Row 26
$ws.Range("H26").Value = 315.86365
$ws.Range("I26").Value = 330.4
$ws.Range("J26").Value = 311.58823
$ws.Range("K26").Value = 991.1999999999999
$ws.Range("L26").Value = 934.76469
$ws.Range("M26").Value = -703.1999999999999
$ws.Range("N26").Value = -1510.76469
# Row 125
$ws.Range("H125").Value = 5799
$ws.Range("I125").Value = 5799
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 17397
$ws.Range("L125").Value = 0
$ws.Range("M125").Value = -12477
$ws.Range("N125").ClearContents()
# Row 136
$ws.Range("H136").Value = 16668172
$ws.Range("I136").Value = 17858234
$ws.Range("K136").Value = 53574702
$ws.Range("M136").Value = -53569602
# Row 137
$ws.Range("H137").Value = 96733.09
$ws.Range("I137").Value = 67516.664
$ws.Range("J137").Value = 169774.17
$ws.Range("K137").Value = 202549.992
$ws.Range("L137").Value = 509322.51
$ws.Range("M137").Value = -197449.992
$ws.Range("N137").Value = -519522.51

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
# Row 31
$ws.Range("H31").Value = 1637
$ws.Range("I31").Value = 1637
$ws.Range("K31").Value = 1637
$ws.Range("M31").Value = -1345
# Row 37
$ws.Range("H37").Value = 1637
$ws.Range("I37").Value = 1637
$ws.Range("K37").Value = 1637
$ws.Range("M37").Value = -1360
# Row 62
$ws.Range("H62").Value = 62232.4
$ws.Range("J62").Value = 70271.25
$ws.Range("L62").Value = 70271.25
$ws.Range("N62").Value = -71643.25
# Row 65
$ws.Range("H65").Value = 62232.4
$ws.Range("J65").Value = 70271.25
$ws.Range("L65").Value = 210813.75
$ws.Range("N65").Value = -217677.75
# Row 113
$ws.Range("H113").Value = 5223.2764
$ws.Range("I113").Value = 2145.4
$ws.Range("K113").Value = 2145.4
$ws.Range("M113").Value = 24.59999999999991
# Row 117
$ws.Range("H117").Value = 56148.332
$ws.Range("J117").Value = 56148.332
$ws.Range("L117").Value = 56148.332
$ws.Range("N117").Value = -63032.332

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 2622
$ws.Range("I22").Value = 499.5
$ws.Range("J22").Value = 4037
$ws.Range("K22").Value = 499.5
$ws.Range("L22").Value = 4037
$ws.Range("M22").Value = -204.5
$ws.Range("N22").Value = -4627
# Row 27
$ws.Range("H27").Value = 2622
$ws.Range("I27").Value = 499.5
$ws.Range("J27").Value = 4037
$ws.Range("K27").Value = 499.5
$ws.Range("L27").Value = 4037
$ws.Range("M27").Value = -392.5
$ws.Range("N27").Value = -4251
# Row 32
$ws.Range("H32").Value = 8802.6
$ws.Range("I32").Value = 8802.6
$ws.Range("K32").Value = 8802.6
$ws.Range("M32").Value = -8485.6
# Row 118
$ws.Range("H118").Value = 56340
$ws.Range("J118").Value = 56340
$ws.Range("L118").Value = 56340
$ws.Range("N118").Value = -59654

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
# Row 39
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()
# Row 64
$ws.Range("H64").Value = 50000
$ws.Range("J64").Value = 50000
$ws.Range("L64").Value = 50000
$ws.Range("N64").Value = -50496
# Row 67
$ws.Range("H67").Value = 50000
$ws.Range("J67").Value = 50000
$ws.Range("L67").Value = 50000
$ws.Range("N67").Value = -51716
# Row 138
$ws.Range("H138").Value = 80000
$ws.Range("J138").Value = 80000
$ws.Range("L138").Value = 80000
$ws.Range("N138").Value = -90280
